# Regenerate save_data: recompute the "K" column (G) values.
# (Per commit message: "regen save_data to use K instead of Strike#,
#  regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 0
    4  = 0
    5  = 2
    6  = 0
    7  = 3
    8  = 0
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    18 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
